$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 'Laptop Lenovo Gaming 15.6'''' IdeaPad 3 15ACH6, FHD IPS, Procesor AMD Ryzen™ 5 5600H (16M Cache, up to 4.2 GHz), 16GB DDR4, 512GB SSD, GeForce RTX 3050 4GB, No OS, Shadow Black'
$ws.Range("B3").Value = 'https://www.pcgarage.ro/notebook-laptop/lenovo/gaming-156-ideapad-3-15ach6-fhd-ips-procesor-amd-ryzen-5-5600h-16m-cache-up-to-42-ghz-16gb-ddr4-512gb-ssd-geforce-rtx-3050-4gb-no-os-shadow-black/'

$ws.Range("A4").Value = 'Laptop Lenovo Gaming 15.6'''' IdeaPad 3 15ACH6, FHD IPS, Procesor AMD Ryzen™ 7 5800H (16M Cache, up to 4.4 GHz), 16GB DDR4, 512GB SSD, GeForce RTX 3050 4GB, No OS, Shadow Black'
$ws.Range("B4").Value = 'https://www.pcgarage.ro/notebook-laptop/lenovo/gaming-156-ideapad-3-15ach6-fhd-ips-procesor-amd-ryzen-7-5800h-16m-cache-up-to-44-ghz-16gb-ddr4-512gb-ssd-geforce-rtx-3050-4gb-no-os-shadow-black/'

$ws.Range("A5").Value = 'Laptop Lenovo Gaming 15.6'''' IdeaPad 3 15ACH6, FHD IPS 120Hz, Procesor AMD Ryzen™ 7 5800H (16M Cache, up to 4.4 GHz), 16GB DDR4, 512GB SSD, GeForce RTX 3050 Ti 4GB, No OS, Shadow Black'
$ws.Range("B5").Value = 'https://www.pcgarage.ro/notebook-laptop/lenovo/gaming-156-ideapad-3-15ach6-fhd-ips-120hz-procesor-amd-ryzen-7-5800h-16m-cache-up-to-44-ghz-16gb-ddr4-512gb-ssd-geforce-rtx-3050-ti-4gb-no-os-shadow-black-1/'

$ws.Range("A6").Value = 'Laptop Lenovo Gaming 15.6'''' IdeaPad 3 15ARH7, FHD IPS 120Hz, Procesor AMD Ryzen™ 5 6600H (16M Cache, up to 4.5 GHz), 16GB DDR5, 512GB SSD, GeForce RTX 3050 Ti 4GB, No OS, Onyx Grey'
$ws.Range("B6").Value = 'https://www.pcgarage.ro/notebook-laptop/lenovo/gaming-156-ideapad-3-15arh7-fhd-ips-120hz-procesor-amd-ryzen-5-6600h-16m-cache-up-to-45-ghz-16gb-ddr5-512gb-ssd-geforce-rtx-3050-ti-4gb-no-os-onyx-grey/'

$ws.Range("A7").Value = 'Laptop Lenovo Gaming 15.6'''' IdeaPad 3 15ARH7, FHD IPS 120Hz, Procesor AMD Ryzen™ 7 6800H (16M Cache, up to 4.7 GHz), 16GB DDR5, 512GB SSD, GeForce RTX 3050 Ti 4GB, No OS, Glacier White'
$ws.Range("B7").Value = 'https://www.pcgarage.ro/notebook-laptop/lenovo/gaming-156-ideapad-3-15arh7-fhd-ips-120hz-procesor-amd-ryzen-7-6800h-16m-cache-up-to-47-ghz-16gb-ddr5-512gb-ssd-geforce-rtx-3050-ti-4gb-no-os-glacier-white/'

$ws.Range("A8").Value = 'Laptop Lenovo Gaming 15.6'''' IdeaPad 3 15ACH6, FHD IPS 120Hz, Procesor AMD Ryzen™ 7 5800H (16M Cache, up to 4.4 GHz), 16GB DDR4, 512GB SSD, GeForce RTX 3050 4GB, No OS, Shadow Black, 4-Zone RGB'
$ws.Range("B8").Value = 'https://www.pcgarage.ro/notebook-laptop/lenovo/gaming-156-ideapad-3-15ach6-fhd-ips-120hz-procesor-amd-ryzen-7-5800h-16m-cache-up-to-44-ghz-16gb-ddr4-512gb-ssd-geforce-rtx-3050-4gb-no-os-shadow-black-4-zone-rgb/'

$ws.Range("A9").Value = 'Laptop Lenovo Gaming 15.6'''' Legion 5 15ARH05, FHD IPS 120Hz, Procesor AMD Ryzen™ 5 4600H (8M Cache, up to 4.0 GHz), 16GB DDR4, 512GB SSD, GeForce GTX 1650 4GB, No OS, Phantom Black'
$ws.Range("B9").Value = 'https://www.pcgarage.ro/notebook-laptop/lenovo/gaming-156-legion-5-15arh05-fhd-ips-120hz-procesor-amd-ryzen-5-4600h-8m-cache-up-to-40-ghz-16gb-ddr4-512gb-ssd-geforce-gtx-1650-4gb-no-os-phantom-black/'

$ws.Range("A10").Value = 'Laptop Lenovo Gaming 15.6'''' IdeaPad 3 15ACH6, FHD IPS 165Hz, Procesor AMD Ryzen™ 5 5600H (16M Cache, up to 4.2 GHz), 16GB DDR4, 512GB SSD, GeForce RTX 3050 Ti 4GB, No OS, Shadow Black, 4-Zone RGB'
$ws.Range("B10").Value = 'https://www.pcgarage.ro/notebook-laptop/lenovo/gaming-156-ideapad-3-15ach6-fhd-ips-165hz-procesor-amd-ryzen-5-5600h-16m-cache-up-to-42-ghz-16gb-ddr4-512gb-ssd-geforce-rtx-3050-ti-4gb-no-os-shadow-black-4-zone-rgb/'

$ws.Range("A11").Value = 'Laptop Lenovo Gaming 15.6'''' Legion 5 15ARH05, FHD IPS 144Hz, Procesor AMD Ryzen™ 5 4600H (8M Cache, up to 4.0 GHz), 16GB DDR4, 512GB SSD, GeForce GTX 1650 4GB, No OS, Phantom Black, 4-Zone RGB'
$ws.Range("B11").Value = 'https://www.pcgarage.ro/notebook-laptop/lenovo/gaming-156-legion-5-15arh05-fhd-ips-144hz-procesor-amd-ryzen-5-4600h-8m-cache-up-to-40-ghz-16gb-ddr4-512gb-ssd-geforce-gtx-1650-4gb-no-os-phantom-black-4-zone-rgb/'

$ws.Range("A12").Value = 'Laptop gaming LENOVO IdeaPad Gaming 3 15IAH7, Intel Core i5-12450H pana la 4.4GHz, 15.6" Full HD, 16GB, SSD 512GB, NVIDIA GeForce RTX 3050 Ti 4GB, Free DOS, Onyx Grey'
$ws.Range("B12").Value = 'https://altex.ro/laptop-gaming-lenovo-ideapad-gaming-3-15iah7-intel-core-i5-12450h-pana-la-4-4ghz-15-6-full-hd-16gb-ssd-512gb-nvidia-geforce-rtx-3050-ti-4gb-free-dos-onyx-grey/cpd/LAP82S900KKRM/'

$ws.Range("A13").Value = 'Laptop Gaming Lenovo IdeaPad 3 15ACH6 cu procesor AMD Ryzen™ 7 5800H pana la 4.40 GHz, 15.6", Full HD, IPS, 16GB, 512GB SSD, NVIDIA GeForce RTX 3050 4GB, No OS, Shadow Black'
$ws.Range("B13").Value = 'https://www.emag.ro/laptop-gaming-lenovo-ideapad-3-15ach6-cu-procesor-amd-ryzentm-7-5800h-pana-la-4-40-ghz-15-6-full-hd-ips-16gb-512gb-ssd-nvidia-geforce-rtx-3050-4gb-no-os-shadow-black-82k201sfrm/pd/DHP2W4MBM/?'

$ws.Range("A14").Value = 'Laptop Gaming Lenovo IdeaPad 3 15ACH6 cu procesor AMD Ryzen™ 7 5800H pana la 4.40 GHz, 15.6'''', Full HD, IPS, 120Hz, 16GB, 512GB SSD, NVIDIA GeForce RTX 3050 Ti 4GB, No OS, Shadow Black'
$ws.Range("B14").Value = 'https://www.emag.ro/laptop-gaming-lenovo-ideapad-3-15ach6-cu-procesor-amd-ryzentm-7-5800h-pana-la-4-40-ghz-15-6-full-hd-ips-120hz-16gb-512gb-ssd-nvidia-geforce-rtx-3050-ti-4gb-no-os-shadow-black-82k201usrm/pd/DN6Y8LMBM/?'

$ws.Range("A15").Value = 'Laptop Gaming Lenovo IdeaPad 3 15ACH6 cu procesor AMD Ryzen™ 5 5600H pana la 4.20 GHz, 15.6", Full HD, IPS, 16GB, 512GB SSD, NVIDIA GeForce RTX 3050 4GB, No OS, Shadow Black'
$ws.Range("B15").Value = 'https://www.emag.ro/laptop-gaming-lenovo-ideapad-3-15ach6-cu-procesor-amd-ryzentm-5-5600h-pana-la-4-20-ghz-15-6-full-hd-ips-16gb-512gb-ssd-nvidia-geforce-rtx-3050-4gb-no-os-shadow-black-82k201sgrm/pd/DF9XMCMBM/?'

$ws.Range("A16").Value = 'Laptop Gaming LENOVO IdeaPad Gaming 3 15ACH6, AMD Ryzen 5 5600H pana la 4.2GHz, 15.6" Full HD, 16GB, SSD 512GB, NVIDIA GeForce GTX 1650, Free DOS, negru'
$ws.Range("B16").Value = 'https://www.emag.ro/laptop-gaming-lenovo-ideapad-gaming-3-15ach6-amd-ryzen-5-5600h-pana-la-4-2ghz-15-6-full-hd-16gb-ssd-512gb-nvidia-geforce-gtx-1650-free-dos-negru-82k2008krm/pd/DCH53JMBM/?'

